# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" worksheet right after "总计" (by copying the
#    existing "2022-Q2" sheet, so it inherits identical styling), then
#    overwrite its data with the Q3 fund-holdings figures.
# 2. Update the "总计" (summary) sheet: shift the existing two data rows
#    down by one and insert the new 2022-Q3 summary row on top.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item("2022-Q2")

# --- 1. Create the new 2022-Q3 sheet as a copy of 2022-Q2 (keeps styles) ---
$q2.Copy($null, $total)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Helper style template: format cells as text ("@") before assignment so
# numeric-looking strings (fund codes with leading zeros, decimal ratios)
# are stored as text instead of being coerced to numbers, then reset the
# style back to Normal so no stray style index is left on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Seed row 6 (new row, doesn't exist yet on the 2022-Q2 template which only
# had 5 rows) by copying row 5's formatting (index-column style s="2"),
# then overwrite with the real values below.
$q3.Range("A5").Copy($q3.Range("A6"))

# Row 2
$q3.Cells.Item(2, 1).Value = 0
Set-TextValue $q3.Cells.Item(2, 2) "550009"
Set-TextValue $q3.Cells.Item(2, 3) "信诚中小盘混合A"
Set-TextValue $q3.Cells.Item(2, 4) "4.09"
Set-TextValue $q3.Cells.Item(2, 5) "92.18"
Set-TextValue $q3.Cells.Item(2, 6) "4.68"
Set-TextValue $q3.Cells.Item(2, 7) "0.1914"
$q3.Cells.Item(2, 8).Value = 5

# Row 3
$q3.Cells.Item(3, 1).Value = 1
Set-TextValue $q3.Cells.Item(3, 2) "004895"
Set-TextValue $q3.Cells.Item(3, 3) "华商鑫安灵活配置混合"
Set-TextValue $q3.Cells.Item(3, 4) "2.11"
Set-TextValue $q3.Cells.Item(3, 5) "92.54"
Set-TextValue $q3.Cells.Item(3, 6) "3.86"
Set-TextValue $q3.Cells.Item(3, 7) "0.0814"
$q3.Cells.Item(3, 8).Value = 9

# Row 4
$q3.Cells.Item(4, 1).Value = 2
Set-TextValue $q3.Cells.Item(4, 2) "005977"
Set-TextValue $q3.Cells.Item(4, 3) "中信保诚至兴灵活配置混合A"
Set-TextValue $q3.Cells.Item(4, 4) "1.21"
Set-TextValue $q3.Cells.Item(4, 5) "92.15"
Set-TextValue $q3.Cells.Item(4, 6) "4.65"
Set-TextValue $q3.Cells.Item(4, 7) "0.0563"
$q3.Cells.Item(4, 8).Value = 4

# Row 5
$q3.Cells.Item(5, 1).Value = 3
Set-TextValue $q3.Cells.Item(5, 2) "016256"
Set-TextValue $q3.Cells.Item(5, 3) "信诚中小盘混合C"
Set-TextValue $q3.Cells.Item(5, 4) "0.45"
Set-TextValue $q3.Cells.Item(5, 5) "92.18"
Set-TextValue $q3.Cells.Item(5, 6) "4.68"
Set-TextValue $q3.Cells.Item(5, 7) "0.0211"
$q3.Cells.Item(5, 8).Value = 5

# Row 6
$q3.Cells.Item(6, 1).Value = 4
Set-TextValue $q3.Cells.Item(6, 2) "005978"
Set-TextValue $q3.Cells.Item(6, 3) "中信保诚至兴灵活配置混合C"
Set-TextValue $q3.Cells.Item(6, 4) "0.38"
Set-TextValue $q3.Cells.Item(6, 5) "92.15"
Set-TextValue $q3.Cells.Item(6, 6) "4.65"
Set-TextValue $q3.Cells.Item(6, 7) "0.0177"
$q3.Cells.Item(6, 8).Value = 4

# --- 2. Update the "总计" sheet summary table -----------------------------
# Final layout:
#   row2: 2022-Q3 (new)
#   row3: 2022-Q2 (was row2)
#   row4: 2021-Q4 (was row3)
# Write bottom-up so we never clobber data we still need to read. Row 4 is
# brand new (the sheet previously only had rows 1-3), so seed its index-
# column (A) style by copying row 3's A-cell formatting first.

$total.Cells.Item(3, 1).Copy($total.Cells.Item(4, 1))
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q4"
$total.Cells.Item(4, 3).Value = 2
$total.Cells.Item(4, 4).Value = 0.19

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 4
$total.Cells.Item(3, 4).Value = 0.68

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.37
